$wb = $excel.ActiveWorkbook

# Update "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# for the zh-cn handback rows, reflecting a newly generated report.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-28 11:41:15"
$wsZhCn.Range("G2").Value = "2016-01-28 11:42:05"
$wsZhCn.Range("D3").Value = "2016-01-28 11:41:15"
$wsZhCn.Range("G3").Value = "2016-01-28 11:42:05"

# Same update for the de-de handback rows.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-28 11:41:28"
$wsDeDe.Range("G2").Value = "2016-01-28 11:42:28"
$wsDeDe.Range("D3").Value = "2016-01-28 11:41:28"
$wsDeDe.Range("G3").Value = "2016-01-28 11:42:28"
